$wb = $excel.ActiveWorkbook

# The new finding is logged on the "hallazgos" sheet, right after the last
# existing row (row 19), so it becomes row 20.
$ws = $wb.Worksheets.Item("hallazgos")

$ws.Range("B20").Value = "correos electronicos en mayusculas"
$ws.Range("C20").Value = "Bajo"
$ws.Range("D20").Value = "Los correos pueden ingresarse en mayusculas y minusculas"
$ws.Range("E20").Value = "se deberá crear una validacion que permita que se ingresen unicamente correos en minusculas"
$ws.Range("F20").Value = "Error de codigo fuente - mantenimiento y validacion"
$ws.Range("G20").Value = "Abierta"

# Match the look of the rows above: wrapped text, title column left
# aligned, row tall enough to show the wrapped description.
$ws.Range("B20").WrapText = $true
$ws.Range("B20").HorizontalAlignment = -4131
$ws.Range("C20:G20").WrapText = $true
$ws.Rows.Item(20).RowHeight = 47.25

# The user ends up on the "hallazgos" tab, scrolled down to the new row,
# with the cell right after the new data selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("G21").Select()

# The other sheet keeps a leftover selection from before the user switched
# tabs to add the new finding.
$ws1 = $wb.Worksheets.Item("Modificacion pre-resolucion")
$ws1.Range("H7").Select()

$wb.Save()
